$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.90202066666667
$ws.Range("H2").Value = 47.706062
$ws.Range("I2").Value = 0.1000095542059894
$ws.Range("J2").Value = 0.10228350994604
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 56.01527666666667
$ws.Range("N2").Value = 168.04583
$ws.Range("O2").Value = 0.4883616647765734
$ws.Range("P2").Value = 0.5237161609491596
$ws.Range("Q2").Value = 890.7560872023845
$ws.Range("R2").Value = 8016.80478482146
$ws.Range("S2").Value = 0.04884083238559997
$ws.Range("T2").Value = 0.05356752715734526

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.90202066666667
$ws.Range("H3").Value = 47.706062
$ws.Range("I3").Value = 0.1000095542059894
$ws.Range("J3").Value = 0.10228350994604
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.149483999999999
$ws.Range("N3").Value = 27.448452
$ws.Range("O3").Value = 0.0797685471532371
$ws.Range("P3").Value = 0.08554331818550501
$ws.Range("Q3").Value = 145.495283657336
$ws.Range("R3").Value = 1309.457552916024
$ws.Range("S3").Value = 0.00797761684045469
$ws.Range("T3").Value = 0.008749670836444366

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.90202066666667
$ws.Range("H4").Value = 47.706062
$ws.Range("I4").Value = 0.1000095542059894
$ws.Range("J4").Value = 0.10228350994604
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.129045
$ws.Range("N4").Value = 36.387135
$ws.Range("O4").Value = 0.1057454494708373
$ws.Range("P4").Value = 0.1134007945935868
$ws.Range("Q4").Value = 192.87632425693
$ws.Range("R4").Value = 1735.88691831237
$ws.Range("S4").Value = 0.01057555526089043
$ws.Range("T4").Value = 0.01159903130170197

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.90202066666667
$ws.Range("H5").Value = 47.706062
$ws.Range("I5").Value = 0.1000095542059894
$ws.Range("J5").Value = 0.10228350994604
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 14.177359
$ws.Range("N5").Value = 42.532077
$ws.Range("O5").Value = 0.1236033999184949
$ws.Range("P5").Value = 0.1325515550349214
$ws.Range("Q5").Value = 225.4486558167527
$ws.Range("R5").Value = 2029.037902350774
$ws.Range("S5").Value = 0.01236152092419331
$ws.Range("T5").Value = 0.01355783829777746

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 15.90202066666667
$ws.Range("H6").Value = 47.706062
$ws.Range("I6").Value = 0.1000095542059894
$ws.Range("J6").Value = 0.10228350994604
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 23.229232
$ws.Range("N6").Value = 46.458464
$ws.Range("O6").Value = 0.2025209386808572
$ws.Range("P6").Value = 0.1447881712368271
$ws.Range("Q6").Value = 369.3917273347947
$ws.Range("R6").Value = 2216.350364008768
$ws.Range("S6").Value = 0.02025402879485105
$ws.Range("T6").Value = 0.01480944235277095

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 24.435136
$ws.Range("H7").Value = 73.305408
$ws.Range("I7").Value = 0.1536752535761215
$ws.Range("J7").Value = 0.1571694269853278
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 56.01527666666667
$ws.Range("N7").Value = 168.04583
$ws.Range("O7").Value = 0.4883616647765734
$ws.Range("P7").Value = 0.5237161609491596
$ws.Range("Q7").Value = 1368.740903427627
$ws.Range("R7").Value = 12318.66813084864
$ws.Range("S7").Value = 0.07504910267139675
$ws.Range("T7").Value = 0.08231216891933511

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 24.435136
$ws.Range("H8").Value = 73.305408
$ws.Range("I8").Value = 0.1536752535761215
$ws.Range("J8").Value = 0.1571694269853278
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 9.149483999999999
$ws.Range("N8").Value = 27.448452
$ws.Range("O8").Value = 0.0797685471532371
$ws.Range("P8").Value = 0.08554331818550501
$ws.Range("Q8").Value = 223.568885869824
$ws.Range("R8").Value = 2012.119972828416
$ws.Range("S8").Value = 0.01225845171117251
$ws.Range("T8").Value = 0.01344479430163939

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 24.435136
$ws.Range("H9").Value = 73.305408
$ws.Range("I9").Value = 0.1536752535761215
$ws.Range("J9").Value = 0.1571694269853278
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 12.129045
$ws.Range("N9").Value = 36.387135
$ws.Range("O9").Value = 0.1057454494708373
$ws.Range("P9").Value = 0.1134007945935868
$ws.Range("Q9").Value = 296.37486412512
$ws.Range("R9").Value = 2667.37377712608
$ws.Range("S9").Value = 0.01625045876195187
$ws.Range("T9").Value = 0.01782313790595489

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 24.435136
$ws.Range("H10").Value = 73.305408
$ws.Range("I10").Value = 0.1536752535761215
$ws.Range("J10").Value = 0.1571694269853278
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.177359
$ws.Range("N10").Value = 42.532077
$ws.Range("O10").Value = 0.1236033999184949
$ws.Range("P10").Value = 0.1325515550349214
$ws.Range("Q10").Value = 346.425695285824
$ws.Range("R10").Value = 3117.831257572416
$ws.Range("S10").Value = 0.01899478382534546
$ws.Range("T10").Value = 0.02083305195085274

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 24.435136
$ws.Range("H11").Value = 73.305408
$ws.Range("I11").Value = 0.1536752535761215
$ws.Range("J11").Value = 0.1571694269853278
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 23.229232
$ws.Range("N11").Value = 46.458464
$ws.Range("O11").Value = 0.2025209386808572
$ws.Range("P11").Value = 0.1447881712368271
$ws.Range("Q11").Value = 567.609443095552
$ws.Range("R11").Value = 3405.656658573312
$ws.Range("S11").Value = 0.03112245660625488
$ws.Range("T11").Value = 0.02275627390754563

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 55.076396
$ws.Range("H12").Value = 165.229188
$ws.Range("I12").Value = 0.3463815024953772
$ws.Range("J12").Value = 0.3542573120827729
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 56.01527666666667
$ws.Range("N12").Value = 168.04583
$ws.Range("O12").Value = 0.4883616647765734
$ws.Range("P12").Value = 0.5237161609491596
$ws.Range("Q12").Value = 3085.119559742893
$ws.Range("R12").Value = 27766.07603768604
$ws.Range("S12").Value = 0.1691594472064532
$ws.Range("T12").Value = 0.1855302794721582

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 55.076396
$ws.Range("H13").Value = 165.229188
$ws.Range("I13").Value = 0.3463815024953772
$ws.Range("J13").Value = 0.3542573120827729
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 9.149483999999999
$ws.Range("N13").Value = 27.448452
$ws.Range("O13").Value = 0.0797685471532371
$ws.Range("P13").Value = 0.08554331818550501
$ws.Range("Q13").Value = 503.9206039796639
$ws.Range("R13").Value = 4535.285435816976
$ws.Range("S13").Value = 0.02763034921481161
$ws.Range("T13").Value = 0.03030434596703839

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 55.076396
$ws.Range("H14").Value = 165.229188
$ws.Range("I14").Value = 0.3463815024953772
$ws.Range("J14").Value = 0.3542573120827729
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 12.129045
$ws.Range("N14").Value = 36.387135
$ws.Range("O14").Value = 0.1057454494708373
$ws.Range("P14").Value = 0.1134007945935868
$ws.Range("Q14").Value = 668.0240855218199
$ws.Range("R14").Value = 6012.21676969638
$ws.Range("S14").Value = 0.03662826766975763
$ws.Range("T14").Value = 0.04017306068077469

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 55.076396
$ws.Range("H15").Value = 165.229188
$ws.Range("I15").Value = 0.3463815024953772
$ws.Range("J15").Value = 0.3542573120827729
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 14.177359
$ws.Range("N15").Value = 42.532077
$ws.Range("O15").Value = 0.1236033999184949
$ws.Range("P15").Value = 0.1325515550349214
$ws.Range("Q15").Value = 780.8378385181639
$ws.Range("R15").Value = 7027.540546663476
$ws.Range("S15").Value = 0.04281393137730526
$ws.Range("T15").Value = 0.04695735759906301

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 55.076396
$ws.Range("H16").Value = 165.229188
$ws.Range("I16").Value = 0.3463815024953772
$ws.Range("J16").Value = 0.3542573120827729
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 23.229232
$ws.Range("N16").Value = 46.458464
$ws.Range("O16").Value = 0.2025209386808572
$ws.Range("P16").Value = 0.1447881712368271
$ws.Range("Q16").Value = 1279.382380407872
$ws.Range("R16").Value = 7676.294282447232
$ws.Range("S16").Value = 0.07014950702704947
$ws.Range("T16").Value = 0.05129226836373862

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 52.98651633333333
$ws.Range("H17").Value = 158.959549
$ws.Range("I17").Value = 0.3332380197777619
$ws.Range("J17").Value = 0.3408149809380521
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 56.01527666666667
$ws.Range("N17").Value = 168.04583
$ws.Range("O17").Value = 0.4883616647765734
$ws.Range("P17").Value = 0.5237161609491596
$ws.Range("Q17").Value = 2968.054372014519
$ws.Range("R17").Value = 26712.48934813067
$ws.Range("S17").Value = 0.1627406741055165
$ws.Range("T17").Value = 0.1784903134108377

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 52.98651633333333
$ws.Range("H18").Value = 158.959549
$ws.Range("I18").Value = 0.3332380197777619
$ws.Range("J18").Value = 0.3408149809380521
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 9.149483999999999
$ws.Range("N18").Value = 27.448452
$ws.Range("O18").Value = 0.0797685471532371
$ws.Range("P18").Value = 0.08554331818550501
$ws.Range("Q18").Value = 484.799283407572
$ws.Range("R18").Value = 4363.193550668148
$ws.Range("S18").Value = 0.02658191269389376
$ws.Range("T18").Value = 0.02915444435677062

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 52.98651633333333
$ws.Range("H19").Value = 158.959549
$ws.Range("I19").Value = 0.3332380197777619
$ws.Range("J19").Value = 0.3408149809380521
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 12.129045
$ws.Range("N19").Value = 36.387135
$ws.Range("O19").Value = 0.1057454494708373
$ws.Range("P19").Value = 0.1134007945935868
$ws.Range("Q19").Value = 642.675841000235
$ws.Range("R19").Value = 5784.082569002116
$ws.Range("S19").Value = 0.03523840418217122
$ws.Range("T19").Value = 0.03864868964777324

# Row 20
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 52.98651633333333
$ws.Range("H20").Value = 158.959549
$ws.Range("I20").Value = 0.3332380197777619
$ws.Range("J20").Value = 0.3408149809380521
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 14.177359
$ws.Range("N20").Value = 42.532077
$ws.Range("O20").Value = 0.1236033999184949
$ws.Range("P20").Value = 0.1325515550349214
$ws.Range("Q20").Value = 751.2088642170304
$ws.Range("R20").Value = 6760.879777953273
$ws.Range("S20").Value = 0.04118935222663803
$ws.Range("T20").Value = 0.04517555570253592

# Row 21
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 52.98651633333333
$ws.Range("H21").Value = 158.959549
$ws.Range("I21").Value = 0.3332380197777619
$ws.Range("J21").Value = 0.3408149809380521
$ws.Range("K21").Value = 2
$ws.Range("M21").Value = 23.229232
$ws.Range("N21").Value = 46.458464
$ws.Range("O21").Value = 0.2025209386808572
$ws.Range("P21").Value = 0.1447881712368271
$ws.Range("Q21").Value = 1230.836080778789
$ws.Range("R21").Value = 7385.016484672737
$ws.Range("S21").Value = 0.06748767656954241
$ws.Range("T21").Value = 0.04934597782013466

# Row 22
$ws.Range("E22").Value = 2
$ws.Range("G22").Value = 10.604946
$ws.Range("H22").Value = 21.209892
$ws.Range("I22").Value = 0.06669566994474985
$ws.Range("J22").Value = 0.04547477004780722
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 56.01527666666667
$ws.Range("N22").Value = 168.04583
$ws.Range("O22").Value = 0.4883616647765734
$ws.Range("P22").Value = 0.5237161609491596
$ws.Range("Q22").Value = 594.0389842250599
$ws.Range("R22").Value = 3564.23390535036
$ws.Range("S22").Value = 0.03257160840760691
$ws.Range("T22").Value = 0.02381587198948343

# Row 23
$ws.Range("E23").Value = 2
$ws.Range("G23").Value = 10.604946
$ws.Range("H23").Value = 21.209892
$ws.Range("I23").Value = 0.06669566994474985
$ws.Range("J23").Value = 0.04547477004780722
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 9.149483999999999
$ws.Range("N23").Value = 27.448452
$ws.Range("O23").Value = 0.0797685471532371
$ws.Range("P23").Value = 0.08554331818550501
$ws.Range("Q23").Value = 97.02978374786399
$ws.Range("R23").Value = 582.178702487184
$ws.Range("S23").Value = 0.005320216692904517
$ws.Range("T23").Value = 0.003890062723612246

# Row 24
$ws.Range("E24").Value = 2
$ws.Range("G24").Value = 10.604946
$ws.Range("H24").Value = 21.209892
$ws.Range("I24").Value = 0.06669566994474985
$ws.Range("J24").Value = 0.04547477004780722
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 12.129045
$ws.Range("N24").Value = 36.387135
$ws.Range("O24").Value = 0.1057454494708373
$ws.Range("P24").Value = 0.1134007945935868
$ws.Range("Q24").Value = 128.62786725657
$ws.Range("R24").Value = 771.7672035394201
$ws.Range("S24").Value = 0.00705276359606619
$ws.Range("T24").Value = 0.005156875057381979

# Row 25
$ws.Range("E25").Value = 2
$ws.Range("G25").Value = 10.604946
$ws.Range("H25").Value = 21.209892
$ws.Range("I25").Value = 0.06669566994474985
$ws.Range("J25").Value = 0.04547477004780722
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 14.177359
$ws.Range("N25").Value = 42.532077
$ws.Range("O25").Value = 0.1236033999184949
$ws.Range("P25").Value = 0.1325515550349214
$ws.Range("Q25").Value = 150.350126617614
$ws.Range("R25").Value = 902.100759705684
$ws.Range("S25").Value = 0.008243811565012858
$ws.Range("T25").Value = 0.006027751484692316

# Row 26
$ws.Range("E26").Value = 2
$ws.Range("G26").Value = 10.604946
$ws.Range("H26").Value = 21.209892
$ws.Range("I26").Value = 0.06669566994474985
$ws.Range("J26").Value = 0.04547477004780722
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 23.229232
$ws.Range("N26").Value = 46.458464
$ws.Range("O26").Value = 0.2025209386808572
$ws.Range("P26").Value = 0.1447881712368271
$ws.Range("Q26").Value = 246.344750981472
$ws.Range("R26").Value = 985.379003925888
$ws.Range("S26").Value = 0.01350726968315938
$ws.Range("T26").Value = 0.006584208792637249
